$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = 15
$ws.Range("C1").Value = 16
$ws.Range("D1").Value = 15
$ws.Range("E1").Value = 16

$ws.Range("B2").Value = 286.48996276264177
$ws.Range("C2").Value = 260.12934017581102
$ws.Range("D2").Value = 287.42060934156262
$ws.Range("E2").Value = 255.41594220805601

$ws.Range("B3").Value = 292.65220380599385
$ws.Range("C3").Value = 249.29294584859031
$ws.Range("D3").Value = 305.96508186043286
$ws.Range("E3").Value = 250.41660847738984

$ws.Range("B1:E3").Select()
